$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.2714889216038224
$ws.Range("B2").Value = -0.3397690405309686
$ws.Range("A3").Value = -0.5046506152826101
$ws.Range("B3").Value = -0.446414115709911
$ws.Range("A4").Value = -0.4967951563606399
$ws.Range("B4").Value = -0.4565314086005962
$ws.Range("A5").Value = -0.1343030677661824
$ws.Range("B5").Value = -0.1802913613953795
$ws.Range("A6").Value = -0.1905517681054592
$ws.Range("B6").Value = -0.1560746376738253
$ws.Range("A7").Value = -0.3242565513960752
$ws.Range("B7").Value = -0.3371946527317504
$ws.Range("A8").Value = -0.5202927461573312
$ws.Range("B8").Value = -0.4519482626887604
$ws.Range("A9").Value = -0.6660051600478981
$ws.Range("B9").Value = -0.5324423675737826
$ws.Range("A10").Value = -0.3870114545086688
$ws.Range("B10").Value = -0.1902536829246043
$ws.Range("A11").Value = -0.2677347412520267
$ws.Range("B11").Value = -0.289408899313559
$ws.Range("A12").Value = -0.1357811179378366
$ws.Range("B12").Value = -0.1441767257263203
$ws.Range("A13").Value = -0.04451853250109153
$ws.Range("B13").Value = -0.05402721325424679
$ws.Range("A14").Value = -0.2232985962213339
$ws.Range("B14").Value = -0.1617285286398429
$ws.Range("A15").Value = -0.0655010880719895
$ws.Range("B15").Value = -0.01043827984171614
$ws.Range("A16").Value = -0.2825864729119449
$ws.Range("B16").Value = -0.1295092751416791
$ws.Range("A17").Value = -0.05852727172672531
$ws.Range("B17").Value = 0.07164748209223724
$ws.Range("A18").Value = 0.1116215181870387
$ws.Range("B18").Value = 0.1076427726448444
$ws.Range("A19").Value = 0.1804747527216866
$ws.Range("B19").Value = 0.1992821640642485
$ws.Range("A20").Value = -0.2776325861956371
$ws.Range("B20").Value = -0.1795743994977774
$ws.Range("A21").Value = 0.0610769799817257
$ws.Range("B21").Value = 0.1170980742962695
$ws.Range("A22").Value = -0.08368337526609629
$ws.Range("B22").Value = 0.06380485484386772
$ws.Range("A23").Value = 0.08126085909716947
$ws.Range("B23").Value = 0.08276718733306503
$ws.Range("A24").Value = 0.8793648853594989
$ws.Range("B24").Value = 0.7433572433828052
$ws.Range("A25").Value = 0.1888685321825936
$ws.Range("B25").Value = 0.1449716609473377
$ws.Range("A26").Value = 0.2134794292027693
$ws.Range("B26").Value = 0.1854128721144535
$ws.Range("A27").Value = 0.1450299496220974
$ws.Range("B27").Value = 0.1431304957956619
$ws.Range("A28").Value = 0.4238204881983129
$ws.Range("B28").Value = 0.2740893809066163
$ws.Range("A29").Value = 0.6967174365804981
$ws.Range("B29").Value = 0.5779327700272949
$ws.Range("A30").Value = 0.2383357926378038
$ws.Range("B30").Value = 0.1994524107542594
$ws.Range("A31").Value = 0.1409859775123565
$ws.Range("B31").Value = 0.1061080477027774
$ws.Range("A32").Value = 0.2039000273958634
$ws.Range("B32").Value = 0.2034803522650228
$ws.Range("A33").Value = 0.08635374808093997
$ws.Range("B33").Value = 0.1021355747555476
$ws.Range("A34").Value = 0.06538158916263581
$ws.Range("B34").Value = 0.04489824213560079
$ws.Range("A35").Value = 0.470348935636911
$ws.Range("B35").Value = 0.3584230454583303
$ws.Range("A36").Value = 0.218918091842702
$ws.Range("B36").Value = 0.1004248712536662
$ws.Range("A37").Value = -0.02345166212809501
$ws.Range("B37").Value = -0.02961433200191043
$ws.Range("A38").Value = 0.2438465863150412
$ws.Range("B38").Value = 0.2063286595340373
$ws.Range("A39").Value = -0.1138876672278304
$ws.Range("B39").Value = -0.1571636705192426
$ws.Range("A40").Value = 0.0787462122475914
$ws.Range("B40").Value = 0.1299504640756246
$ws.Range("A41").Value = 0.01190084292612895
$ws.Range("B41").Value = -0.1020176820699642
$ws.Range("A42").Value = 0.3306898555637756
$ws.Range("B42").Value = 0.3177028113401377
$ws.Range("A43").Value = -0.008900747196598885
$ws.Range("B43").Value = 0.02457878542998314
$ws.Range("A44").Value = 0.08762448441798998
$ws.Range("B44").Value = 0.006880123792583711
$ws.Range("A45").Value = -0.1193520390177986
$ws.Range("B45").Value = -0.08089627767262299
$ws.Range("A46").Value = -0.1761564515494458
$ws.Range("B46").Value = -0.1589291480285107
$ws.Range("A47").Value = -0.1787085014596541
$ws.Range("B47").Value = -0.1580723057225032
$ws.Range("A48").Value = -0.2274787900245934
$ws.Range("B48").Value = -0.1977470538024043
$ws.Range("A49").Value = -0.2261173547094197
$ws.Range("B49").Value = -0.1918145497072672
$ws.Range("A50").Value = -0.07576855223009549
$ws.Range("B50").Value = -0.08459245785905996
$ws.Range("A51").Value = -0.2516436965521065
$ws.Range("B51").Value = -0.2274791057623393
$ws.Range("A52").Value = -0.2516436965521065
$ws.Range("B52").Value = -0.2274791057623393
$ws.Range("A53").Value = -0.2115869172361382
$ws.Range("B53").Value = -0.1864422125477097
$ws.Range("A54").Value = -0.2030205474110467
$ws.Range("B54").Value = -0.1735016429854447
$ws.Range("A55").Value = -0.1698698248759659
$ws.Range("B55").Value = -0.1427455500046938
$ws.Range("A56").Value = -0.09628301818186084
$ws.Range("B56").Value = -0.08553992600288009
$ws.Range("A57").Value = -0.1905854495412856
$ws.Range("B57").Value = -0.1423186524008358
$ws.Range("A58").Value = -0.1272014980566598
$ws.Range("B58").Value = -0.1645691247519662
$ws.Range("A59").Value = -0.2090670847908934
$ws.Range("B59").Value = -0.2153037535768792
$ws.Range("A60").Value = -0.2485626281657644
$ws.Range("B60").Value = -0.2376797652111406
$ws.Range("A61").Value = -0.279213867352844
$ws.Range("B61").Value = -0.2209887020818051
$ws.Range("A62").Value = -0.1511940731555665
$ws.Range("B62").Value = -0.08489015703895608
$ws.Range("A63").Value = -0.5138123228454438
$ws.Range("B63").Value = -0.4684567990295513
$ws.Range("A64").Value = -0.3196448844002049
$ws.Range("B64").Value = -0.2932362443483825
$ws.Range("A65").Value = -0.1770104367733414
$ws.Range("B65").Value = -0.1751837628722738
$ws.Range("A66").Value = -0.07477732495492206
$ws.Range("B66").Value = -0.05613786915979644
$ws.Range("A67").Value = 0.07780695627103607
$ws.Range("B67").Value = 0.05692010331208052
